# v0.1.7 - Renamed RI columns to Total Pagado/Pendiente
# - Replaced all RI Pagado references with Total Pagado
# - Replaced all RI Pendiente references with Total Pendiente
# - Updated KPI cards, table headers, and chart titles
# - Updated Excel data with new column names

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pagado_pendiente")

# Rename the "Recibo_Inicial_*" headers (and the backing Tabla2 columns)
# to "Total_Pagado" / "Total_Pendiente". Writing straight to the header
# cells keeps the table definition (xl/tables/table1.xml) in sync, the
# same way Excel itself would when you retype a table header.
$ws.Range("D2").Value = "Total_Pagado"
$ws.Range("F2").Value = "Total_Pendiente"

# Refresh the underlying figures for the renamed columns.
$ws.Range("D3").Value = 104105.63
$ws.Range("F3").Value = 225651.58
$ws.Range("D4").Value = 74359.21
$ws.Range("F4").Value = 108385.48
$ws.Range("D5").Value = 71251.41
$ws.Range("F5").Value = 37782.11
$ws.Range("D6").Value = 52982.99
$ws.Range("F6").Value = 21966.88
$ws.Range("D7").Value = 40909.81
$ws.Range("F7").Value = 204839.8
$ws.Range("D8").Value = 27079.33
$ws.Range("F8").Value = 801637.94
$ws.Range("D9").Value = 26005.4
$ws.Range("F9").Value = 99845.58
$ws.Range("D10").Value = 22054.4
$ws.Range("F10").Value = 75875.25
$ws.Range("D11").Value = 20077.4
$ws.Range("F11").Value = 42810.25
$ws.Range("D12").Value = 14887.38
$ws.Range("F12").Value = 14887.38
$ws.Range("D13").Value = 11945.29
$ws.Range("F13").Value = 66315.25
$ws.Range("D14").Value = 11532.26
$ws.Range("F14").Value = 58958.98
$ws.Range("D15").Value = 9688.87
$ws.Range("F15").Value = 24855.21
$ws.Range("D16").Value = 9281.48
$ws.Range("F16").Value = 34668.69
$ws.Range("D17").Value = 8821.32
$ws.Range("F17").Value = 11272.92
$ws.Range("D18").Value = 7991.52
$ws.Range("F18").Value = 55306.96
$ws.Range("D19").Value = 7460.84
$ws.Range("F19").Value = 15352.07
$ws.Range("D20").Value = 6716.14
$ws.Range("F20").Value = 43873.2
$ws.Range("D21").Value = 6639.64
$ws.Range("F21").Value = 8336.92
$ws.Range("D22").Value = 6196.74
$ws.Range("D23").Value = 5895.02
$ws.Range("F23").Value = 88859
$ws.Range("D24").Value = 5745.69
$ws.Range("F24").Value = 22496
$ws.Range("D25").Value = 5138.51
$ws.Range("F25").Value = 10277.03
$ws.Range("D26").Value = 4663.49
$ws.Range("F26").Value = 3167.47
$ws.Range("D27").Value = 4487.59
$ws.Range("F27").Value = 38209.88
$ws.Range("D28").Value = 4192.72
$ws.Range("F28").Value = 23946.42
$ws.Range("D29").Value = 4095.39
$ws.Range("F29").Value = 63114.42
$ws.Range("D30").Value = 3993.98
$ws.Range("F30").Value = 26539.81
$ws.Range("D31").Value = 3904.72
$ws.Range("F31").Value = 49420.01
$ws.Range("D32").Value = 3151.08
$ws.Range("D33").Value = 2852.69
$ws.Range("F33").Value = 1575.88
$ws.Range("D34").Value = 2833.03
$ws.Range("F34").Value = 4506.84
$ws.Range("D35").Value = 2133.91
$ws.Range("D36").Value = 1721.93
$ws.Range("F36").Value = 22676.73
$ws.Range("D37").Value = 1712.03
$ws.Range("D38").Value = 1345.8
$ws.Range("F38").Value = 4779.23
$ws.Range("D39").Value = 1297.97
$ws.Range("F39").Value = 32930.35
$ws.Range("D40").Value = 1269.82
$ws.Range("F40").Value = 49303.64
$ws.Range("D41").Value = 651.19
$ws.Range("F41").Value = 46859.36
$ws.Range("F42").Value = 18020.83
$ws.Range("F43").Value = 127958.49
$ws.Range("F44").Value = 10092.96
$ws.Range("F45").Value = 15118.28
$ws.Range("F46").Value = 12739.06

# Move the selection cursor as it was left after the edit.
$ws.Range("H13").Select()
